$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1) URL property: matchsource -> matchsync
$ws.Range("B2").Value = "http://fhir.nmdp.org/ig/matchsync/ValueSet/nmdp-transplant-timeline-codes"

# 2) Experimental property: now has a literal text value "true" (not boolean TRUE).
#    A bare "true"/"false" typed into a cell is auto-coerced to a Boolean, so
#    force text interpretation with a quote-prefix, then copy the (unchanged)
#    plain format from the cell above back onto it, so the cell keeps its
#    original border/alignment style instead of the quote-prefix style.
$ws.Cells.Item(7, 2).Value = "'true"
$ws.Cells.Item(7, 1).Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3) Date property: refreshed publication date
$ws.Range("B8").Value = "2024-02-19T18:37:26-06:00"
